$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("L2").Value = 0.7808705382933534
$ws.Range("Q2").Value = 0.02351246133036713

# Row 3
$ws.Range("L3").Value = 0.7808705382933534
$ws.Range("Q3").Value = 0.02351246133036713

# Row 4
$ws.Range("L4").Value = 0.7832122299475502
$ws.Range("Q4").Value = -0.01081937260331701
